# Generate Report for handback
# Adds a new handback entry (35915d0c-9835-4a7c-8dc4-ecd480c1bdac) as row 4
# to the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$guid = "35915d0c-9835-4a7c-8dc4-ecd480c1bdac"
$mdName = "$guid.md"

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(4, 1).Value = $mdName
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$mdName", "", "", $mdName)
$wsOverview.Cells.Item(4, 2).Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item(4, 3).Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhXlfName = "$guid.f1eb4b1c6ce5e6fc8218d35db879bc037c7dd809.zh-cn.xlf"

$wsZhCn.Cells.Item(4, 1).Value = $mdName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$mdName", "", "", $mdName)
$wsZhCn.Cells.Item(4, 2).Value = "Handed back: in sync with en-US"
$wsZhCn.Cells.Item(4, 3).Value = $zhXlfName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlfName", "", "", $zhXlfName)
$wsZhCn.Cells.Item(4, 4).Value = "2016-01-18 06:29:43"
$wsZhCn.Cells.Item(4, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(4, 5).Value = $mdName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0000000000000000000000000000000000000000/e2e/$mdName", "", "", $mdName)
$wsZhCn.Cells.Item(4, 6).Value = $zhXlfName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000000/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlfName", "", "", $zhXlfName)
$wsZhCn.Cells.Item(4, 7).Value = "2016-01-18 06:30:25"
$wsZhCn.Cells.Item(4, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(4, 8).Value = "Include"

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deXlfName = "$guid.f1eb4b1c6ce5e6fc8218d35db879bc037c7dd809.de-de.xlf"

$wsDeDe.Cells.Item(4, 1).Value = $mdName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$mdName", "", "", $mdName)
$wsDeDe.Cells.Item(4, 2).Value = "Handed back: in sync with en-US"
$wsDeDe.Cells.Item(4, 3).Value = $deXlfName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlfName", "", "", $deXlfName)
$wsDeDe.Cells.Item(4, 4).Value = "2016-01-18 06:29:53"
$wsDeDe.Cells.Item(4, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(4, 5).Value = $mdName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/0000000000000000000000000000000000000000/e2e/$mdName", "", "", $mdName)
$wsDeDe.Cells.Item(4, 6).Value = $deXlfName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0000000000000000000000000000000000000000/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlfName", "", "", $deXlfName)
$wsDeDe.Cells.Item(4, 7).Value = "2016-01-18 06:30:41"
$wsDeDe.Cells.Item(4, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(4, 8).Value = "Include"
